$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; this shifts rows 26..109 down to 27..110
# and mirrors Excel's native "insert row" behaviour (formatting copied down).
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new data record.
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44998
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 100112040
$ws.Range("G26").Value = "Cilantro"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 1500
$ws.Range("L26").Value = 2000
$ws.Range("M26").Value = 1750
$ws.Range("N26").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 875
$ws.Range("Q26").Value = 2
$ws.Range("R26").Value = "Hortaliza"
